# Test script: insert row for AllCare To You and append row for Western Health Advantage
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (shifts existing rows 5.. down by one)
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "AllCare To You"
$ws.Cells.Item(5, 2).Value = 1750
$ws.Cells.Item(5, 3).Value = "Cozeva Support"
$ws.Cells.Item(5, 4).Value = 99999

# Append new row at the end
$ws.Cells.Item(132, 1).Value = "Western Health Advantage"
$ws.Cells.Item(132, 2).Value = 7200
$ws.Cells.Item(132, 3).Value = "Cozeva Support"
$ws.Cells.Item(132, 4).Value = 99999
